# Weekly price update for Hortaliza, Vega Modelo de Temuco - Acelga
# A new weekly record is inserted at row 262, pushing the existing
# rows 262-278 down to rows 263-279.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 262 (shifts rows 262:278 down to 263:279)
$ws.Rows(262).Insert()

# Populate the newly inserted row 262 with the new weekly record
$ws.Range("A262").Value = 10
$ws.Range("B262").Value = "Vega Modelo de Temuco"
$ws.Range("C262").Value = "La Araucanía"
$ws.Range("D262").Value = 44610
$ws.Range("E262").Value = 9
$ws.Range("F262").Value = 100112009
$ws.Range("G262").Value = "Acelga"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 40
$ws.Range("K262").Value = 8000
$ws.Range("L262").Value = 8000
$ws.Range("M262").Value = 8000
$ws.Range("N262").Value = "$/docena de atados (12 kilos)"
$ws.Range("O262").Value = "Provincia de Cautín"
$ws.Range("P262").Value = 667
$ws.Range("Q262").Value = 12
$ws.Range("R262").Value = "Hortaliza"
